$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bus")

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "J2" "-19.908076343375843"
Set-TextValue "K2" "-57.80509971944206"
Set-TextValue "J3" "19.908076343375843"
Set-TextValue "K3" "57.80509971944206"
